# Update the marker_table genotype column: collapse homozygous
# double-letter genotype calls (e.g. "GG") down to a single letter ("G").
$wb = $excel.ActiveWorkbook
$markerSheet = $wb.Worksheets.Item("marker_table")

$genotypeUpdates = @{
    2  = "G"
    3  = "C"
    4  = "T"
    5  = "G"
    6  = "G"
    7  = "G"
    8  = "G"
    9  = "T"
    10 = "G"
    12 = "G"
    13 = "A"
    14 = "C"
    15 = "A"
    16 = "T"
    17 = "G"
    18 = "G"
}

foreach ($row in $genotypeUpdates.Keys) {
    $markerSheet.Range("G$row").Value = $genotypeUpdates[$row]
}

# Fill in the computed diplotype result for the sample on genotype_result.
$resultSheet = $wb.Worksheets.Item("genotype_result")
$resultSheet.Range("B2").Value = "*2/*5"
